# Remove the trailing site-chrome paragraphs that were scraped along with
# the bibliography ("Ver no Jupiter ..." and the "(c) 2020 ..." footer
# line), together with the blank paragraph that separated them from the
# last bibliography entry. The blank paragraph right before the
# page-break paragraph at the very end of the document is left
# untouched.

$d = $word.ActiveDocument

$i = 0
$startIndex = -1
$endIndex = -1
ForEach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startIndex = $i - 1
        $endIndex = $i
    } elseif ($t -like "*luizeleno@usp.br*") {
        $endIndex = $i
    }
}

$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($endIndex)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
